# ---------------------------------------------------------------------------
# Commit: "test updates because of change to tax calc"
#
# The DISCK DCF model's tax-calculation logic changed upstream, which shifts
# every downstream projection figure (tax, cash build-up, diluted share count,
# implied price, FCFE/FCFF, dividends, equity/EV/WACC, DDM and per-share value)
# on the 'raw data' sheet, plus the three cached headline figures on 'report'
# (the as-of date used for the projection, the DCF cash-generated figure, and
# the DDM intrinsic value per share). This script pokes in the newly
# recalculated literal values for every cell the model re-generated.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 'report' sheet: cached headline summary values -------------------------
$wsReport = $wb.Worksheets.Item("report")
$wsReport.Range("B3").Value = 45621.0  # Date of preparation
$wsReport.Range("B11").Value = 146.8567227010758  # Cash Generated
$wsReport.Range("B17").Value = 128.6959462351709  # IV per share DDM model

# --- 'raw data' sheet: recalculated 11-year projection table -----------------
$wsRaw = $wb.Worksheets.Item("raw data")
# row 7: tax
$wsRaw.Range("D7").Value = 0.81345984825
$wsRaw.Range("E7").Value = 0.9292499999999998
$wsRaw.Range("F7").Value = 0.9683625
# row 9: debt
$wsRaw.Range("C9").Value = 43.482795
# row 10: interest
$wsRaw.Range("D10").Value = 2.826381675
# row 11: cash
$wsRaw.Range("C11").Value = -0.0000000000000004440892098500626
$wsRaw.Range("D11").Value = 0.2773634767499966
$wsRaw.Range("E11").Value = 11.92311347675
$wsRaw.Range("F11").Value = 26.01600097674999
$wsRaw.Range("G11").Value = 43.1530139160036
$wsRaw.Range("H11").Value = 59.45745229451082
$wsRaw.Range("I11").Value = 77.39334111227166
$wsRaw.Range("J11").Value = 97.8171678692861
$wsRaw.Range("K11").Value = 120.6307450655541
$wsRaw.Range("L11").Value = 146.8567227010758
# row 14: shares
$wsRaw.Range("D14").Value = 2.291472942632959
$wsRaw.Range("E14").Value = 1.992339438283776
$wsRaw.Range("F14").Value = 1.707127266738211
$wsRaw.Range("G14").Value = 1.448094808099672
$wsRaw.Range("H14").Value = 1.251937235886192
$wsRaw.Range("I14").Value = 1.081599197054031
$wsRaw.Range("J14").Value = 0.9310767383956051
$wsRaw.Range("K14").Value = 0.8014168273461849
$wsRaw.Range("L14").Value = 0.8014168273461849
# row 15: price
$wsRaw.Range("C15").Value = 32.5274552299937
$wsRaw.Range("D15").Value = 38.93161357948631
$wsRaw.Range("E15").Value = 49.41194277800501
$wsRaw.Range("F15").Value = 66.15778203752855
$wsRaw.Range("G15").Value = 83.11908734658975
$wsRaw.Range("H15").Value = 105.2958513596226
$wsRaw.Range("I15").Value = 135.6862420335648
$wsRaw.Range("J15").Value = 175.9493509722724
$wsRaw.Range("K15").Value = 230.7867578932841
$wsRaw.Range("L15").Value = 230.7867578932841
# row 17: buybacks
$wsRaw.Range("C17").Value = -0.0000000000000004440892098500626
$wsRaw.Range("D17").Value = 0.2773634767499971
# row 19: fcfe
$wsRaw.Range("C19").Value = -0.0000000000000004440892098500626
$wsRaw.Range("D19").Value = 0.2773634767499971
# row 20: fcff
$wsRaw.Range("C20").Value = 10.872
# row 21: fcf
$wsRaw.Range("C21").Value = 8.217205
$wsRaw.Range("D21").Value = 8.76015847675
# row 22: dividend
$wsRaw.Range("L22").Value = 32.72451580829226
# row 23: income_pretax
$wsRaw.Range("D23").Value = 3.873618325
# row 24: dDebt
$wsRaw.Range("C24").Value = -8.217205
$wsRaw.Range("D24").Value = -8.482795000000003
# row 25: tax_cash
$wsRaw.Range("C25").Value = 0.722295
$wsRaw.Range("D25").Value = 0.81345984825
$wsRaw.Range("E25").Value = 0.9292499999999998
$wsRaw.Range("F25").Value = 0.9683625
# row 26: income_taxable
$wsRaw.Range("D26").Value = 3.873618325
# row 29: equity
$wsRaw.Range("B29").Value = 169.994742848637
$wsRaw.Range("C29").Value = 186.9942171335007
# row 30: EV
$wsRaw.Range("B30").Value = 221.694742848637
$wsRaw.Range("C30").Value = 230.4770121335007
# row 31: wacc
$wsRaw.Range("B31").Value = 0.08865464752261958
$wsRaw.Range("C31").Value = 0.09082147951690445
# row 32: firm
$wsRaw.Range("B32").Value = 222.6781438463931
$wsRaw.Range("C32").Value = 231.5475962000862
# row 33: DDM
$wsRaw.Range("B33").Value = 128.6959462351709
$wsRaw.Range("C33").Value = 141.5655408586879
$wsRaw.Range("D33").Value = 155.7220949445567
$wsRaw.Range("E33").Value = 171.2943044390124
$wsRaw.Range("F33").Value = 188.4237348829137
$wsRaw.Range("G33").Value = 207.266108371205
$wsRaw.Range("H33").Value = 227.9927192083256
$wsRaw.Range("I33").Value = 250.7919911291581
$wsRaw.Range("J33").Value = 275.871190242074
$wsRaw.Range("K33").Value = 303.4583092662814
$wsRaw.Range("L33").Value = 301.0796243846173
# row 34: value_per_share
$wsRaw.Range("B34").Value = 73.91075776027695
$wsRaw.Range("C34").Value = 81.30183353630466
# row 35: value_per_share_DDM
$wsRaw.Range("B35").Value = 128.6959462351709
$wsRaw.Range("C35").Value = 141.5655408586879
$wsRaw.Range("D35").Value = 155.7220949445567
$wsRaw.Range("E35").Value = 171.2943044390124
$wsRaw.Range("F35").Value = 188.4237348829137
$wsRaw.Range("G35").Value = 207.266108371205
$wsRaw.Range("H35").Value = 227.9927192083256
$wsRaw.Range("I35").Value = 250.7919911291581
$wsRaw.Range("J35").Value = 275.871190242074
$wsRaw.Range("K35").Value = 303.4583092662814
$wsRaw.Range("L35").Value = 301.0796243846173
